$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E6").Value = "Four digit calendar year (ex: 2014, 2016)"
$ws.Range("F11").Value = "Total tonnage of contaminated material (from city data); usually refers to hazardous waste."
$ws.Range("F12").Value = "Amount recycled in tons for the given year in each city (cumulative residential totals; hand calculated for Portland and Los Angeles)."

$ws.Columns.Item(5).ColumnWidth = 18.166666666666668

$ws.Rows.Item(6).RowHeight = 33.6
$ws.Rows.Item(12).RowHeight = 43.2
$ws.Rows.Item(13).RowHeight = 34.200000000000003
$ws.Rows.Item(15).RowHeight = 82.2
$ws.Rows.Item(17).RowHeight = 36

$ws.Range("A4").Select()

$wb.Save()
